$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.509.13"
$ws.Range("E2").Value = "  -2.28%  "
$ws.Range("D3").Value = "1.959.28"
$ws.Range("E3").Value = "  -3.92%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'250.77"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").Value = "'0.603"
$ws.Range("E6").Value = "  -3.45%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'53.96"
$ws.Range("E8").Value = "  -8.72%  "
$ws.Range("D9").Value = "'0.368"
$ws.Range("E9").Value = "  -6.39%  "
$ws.Range("D10").Value = "'0.0748"
$ws.Range("E10").Value = "  -7.26%  "
$ws.Range("E11").Value = "  -3.85%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.251.75"
$ws.Range("E12").Value = "  -3.86%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'13.93"
$ws.Range("E13").Value = "  -8.48%  "
$ws.Range("D14").Value = "'21.03"
$ws.Range("E14").Value = "  -4.54%  "
$ws.Range("D15").Value = "'0.759"
$ws.Range("E15").Value = "  -10.83%  "
$ws.Range("D16").Value = "'5.08"
$ws.Range("E16").Value = "  -6.81%  "
$ws.Range("D17").Value = "1.963.52"
$ws.Range("E17").Value = "  -3.57%  "
$ws.Range("D18").Value = "36.393.19"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("D19").Value = "'68.57"
$ws.Range("E19").Value = "  -2.72%  "
$ws.Range("E20").Value = "  -5.65%  "
$ws.Range("D21").Value = "'229.04"
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("E22").Value = "  -5.67%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'2.48"
$ws.Range("E24").Value = "  -2.73%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'162.73"
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("D27").Value = "'8.70"
$ws.Range("E27").Value = "  -7.43%  "
$ws.Range("D28").Value = "'18.96"
$ws.Range("E28").Value = "  -4.85%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.121"
$ws.Range("E29").Value = "  -12.25%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.31"
$ws.Range("E30").Value = "  -4.77%  "
$ws.Range("E31").Value = "  -3.74%  "
$ws.Range("D32").Value = "'4.44"
$ws.Range("E32").Value = "  -7.12%  "
$ws.Range("D33").Value = "'0.0616"
$ws.Range("E33").Value = "  -9.76%  "
$ws.Range("D34").Value = "'4.26"
$ws.Range("E34").Value = "  -5.94%  "
$ws.Range("D35").Value = "'2.33"
$ws.Range("E35").Value = "  -9.09%  "
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "'3.33"
$ws.Range("E38").Value = "  -7.03%  "
$ws.Range("D39").Value = "'5.29"
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("D41").Value = "1.439.74"
$ws.Range("E41").Value = "  +3.67%  "
$ws.Range("D42").Value = "'1.14"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("D43").Value = "'0.0898"
$ws.Range("E43").Value = "  -8.48%  "
$ws.Range("E44").Value = "  -6.05%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'86.77"
$ws.Range("E45").Value = "  -5.48%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'15.21"
$ws.Range("E46").Value = "  -8.69%  "
$ws.Range("E47").Value = "  -6.34%  "
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").Value = "'6.73"
$ws.Range("E49").Value = "  -10.61%  "
$ws.Range("D50").Value = "2.144.55"
$ws.Range("E50").Value = "  -4.00%  "
$ws.Range("D51").Value = "'1.88"
$ws.Range("E51").Value = "  -10.95%  "
